# Add slide titles for all slides (and body content for the "Free RCS
# hosts" slide). Ten new slides are appended after the existing title
# slide, each using the "Title and Content" layout (ppLayoutText = 2),
# matching the target deck's slide2.xml .. slide11.xml.

$p = $ppt.ActivePresentation

$slides = @(
    @{ Title = "About versions";               Body = $null },
    @{ Title = "Working together";              Body = $null },
    @{ Title = "Version control systems";        Body = $null },
    @{ Title = "CVS";                            Body = $null },
    @{ Title = "SVN";                            Body = $null },
    @{ Title = "Distributed version control";    Body = $null },
    @{ Title = "GIT";                            Body = $null },
    @{ Title = "Mercurial";                      Body = $null },
    @{ Title = "Free RCS hosts";                 Body = "SourceForge`nGitHub`nBitBucket" },
    @{ Title = "Exercise: basic git usage";      Body = $null }
)

for ($i = 0; $i -lt $slides.Length; $i++) {
    $slideIndex = $i + 2
    $info = $slides[$i]

    $s = $p.Slides.Add($slideIndex, 2)
    $s.Shapes.Item(1).TextFrame.TextRange.Text = $info.Title

    if ($info.Body) {
        $s.Shapes.Item(2).TextFrame.TextRange.Text = $info.Body
    }
}
